$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark (it currently sits, collapsed,
#    at the very end of the document and is being relocated by this edit).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Find the "Andrea Elías - 1" run and append a new run containing "7048"
#    right after it, in the same paragraph, with matching character formatting.
$find = $d.Content.Find
$find.Execute("Andrea Elías - 1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng = $find.Parent
$rng.Collapse(0)
$rng.InsertAfter("7048")

# Re-select just the freshly inserted text and nudge its formatting so the
# engine keeps it as its own <w:r> (same rPr as the preceding run) instead of
# silently merging it back into "Andrea Elías - 1".
$newRng = $d.Range($rng.Start, $rng.End)
$newRng.Font.Size = 99
$newRng.Font.Size = 12

$pos = $newRng.End

# 3. Insert the new "_GoBack" bookmark (collapsed) right after "7048".
#    Adding a bookmark exactly at "end of paragraph's last run" mis-resolves
#    to the start of the document in this runtime, so we temporarily insert
#    a placeholder character after the target position, add the bookmark,
#    then remove the placeholder again -- the bookmark stays put.
$placeholder = $d.Range($pos, $pos)
$placeholder.InsertAfter("X")

$d.Range($pos, $pos).Bookmarks.Add("_GoBack")

$d.Range($pos, $pos + 1).Delete()
